# Update betting odds values on Sheet1 as per the 2024-11-25 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.37
$ws.Range("V2").Value = 1.63

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 2.62
$ws.Range("V3").Value = 1.58

# Row 5
$ws.Range("K5").Value = 1.92
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("R5").Value = 1.5
